# ============================================================================
# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board pricing/profit columns (H:N) for the
# Leve profit tables on every job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Columns:
#   H currentAveragePrice    I currentAveragePriceNQ  J currentAveragePriceHQ
#   K LevePriceNQ            L LevePriceHQ            M LeveProfitNQ
#   N LeveProfitHQ
# Only the numeric market-data cells change; no rows/columns are added or
# removed and no formulas/styles are touched.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 18
$ws.Cells.Item(18, 8).Value = 393.33334
$ws.Cells.Item(18, 9).Value = 393.33334
$ws.Cells.Item(18, 11).Value = 393.33334
$ws.Cells.Item(18, 13).Value = -109.33334
# Row 19
$ws.Cells.Item(19, 8).Value = 1816.1
$ws.Cells.Item(19, 10).Value = 2551.2222
$ws.Cells.Item(19, 12).Value = 2551.2222
$ws.Cells.Item(19, 14).Value = -2901.2222
# Row 95
$ws.Cells.Item(95, 8).Value = 749749
$ws.Cells.Item(95, 10).Value = 749749
$ws.Cells.Item(95, 12).Value = 749749
$ws.Cells.Item(95, 14).Value = -755241
# Row 113
$ws.Cells.Item(113, 8).Value = 6649.1763
$ws.Cells.Item(113, 9).Value = 4755.625
$ws.Cells.Item(113, 10).Value = 8332.333
$ws.Cells.Item(113, 11).Value = 4755.625
$ws.Cells.Item(113, 12).Value = 8332.333
$ws.Cells.Item(113, 13).Value = -1501.625
$ws.Cells.Item(113, 14).Value = -14840.333
# Row 129
$ws.Cells.Item(129, 8).Value = 1608.2
$ws.Cells.Item(129, 10).Value = 8000
$ws.Cells.Item(129, 12).Value = 24000
$ws.Cells.Item(129, 14).Value = -34000
# Row 135
$ws.Cells.Item(135, 8).Value = 909.0769
$ws.Cells.Item(135, 9).Value = 729.381
$ws.Cells.Item(135, 10).Value = 1663.8
$ws.Cells.Item(135, 11).Value = 6564.429
$ws.Cells.Item(135, 12).Value = 14974.2
$ws.Cells.Item(135, 13).Value = -4029.429
$ws.Cells.Item(135, 14).Value = -20044.2
# Row 141
$ws.Cells.Item(141, 8).Value = 1773.95
$ws.Cells.Item(141, 9).Value = 1687.6471
$ws.Cells.Item(141, 11).Value = 5062.9413
$ws.Cells.Item(141, 13).Value = 117.0587000000005

# ---------------------------------------------------------------------------
# ARM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Cells.Item(32, 8).Value = 3887.5647
$ws.Cells.Item(32, 9).Value = 2175.6216
$ws.Cells.Item(32, 11).Value = 2175.6216
$ws.Cells.Item(32, 13).Value = -1888.6216
# Row 96
$ws.Cells.Item(96, 8).Value = 64911.5
$ws.Cells.Item(96, 10).Value = 64911.5
$ws.Cells.Item(96, 12).Value = 64911.5
$ws.Cells.Item(96, 14).Value = -70403.5
# Row 97
$ws.Cells.Item(97, 8).Value = 1094.7241
$ws.Cells.Item(97, 9).Value = 999.12
$ws.Cells.Item(97, 10).Value = 1692.25
$ws.Cells.Item(97, 11).Value = 999.12
$ws.Cells.Item(97, 12).Value = 1692.25
$ws.Cells.Item(97, 13).Value = -503.12
$ws.Cells.Item(97, 14).Value = -2684.25

# ---------------------------------------------------------------------------
# BSM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Cells.Item(20, 8).Value = 2050.2273
$ws.Cells.Item(20, 9).Value = 1901.1333
$ws.Cells.Item(20, 11).Value = 1901.1333
$ws.Cells.Item(20, 13).Value = -1654.1333
# Row 21
$ws.Cells.Item(21, 8).Value = 597597
$ws.Cells.Item(21, 10).Value = 597597
$ws.Cells.Item(21, 12).Value = 597597
$ws.Cells.Item(21, 14).Value = -598069
# Row 74
$ws.Cells.Item(74, 8).Value = 44859.5
$ws.Cells.Item(74, 9).Value = 29700
$ws.Cells.Item(74, 10).Value = 47891.4
$ws.Cells.Item(74, 11).Value = 29700
$ws.Cells.Item(74, 12).Value = 47891.4
$ws.Cells.Item(74, 13).Value = -28764
$ws.Cells.Item(74, 14).Value = -49763.4
# Row 77
$ws.Cells.Item(77, 8).Value = 44859.5
$ws.Cells.Item(77, 9).Value = 29700
$ws.Cells.Item(77, 10).Value = 47891.4
$ws.Cells.Item(77, 11).Value = 89100
$ws.Cells.Item(77, 12).Value = 143674.2
$ws.Cells.Item(77, 13).Value = -84420
$ws.Cells.Item(77, 14).Value = -153034.2
# Row 87
$ws.Cells.Item(87, 8).Value = 84357
$ws.Cells.Item(87, 9).Value = 86199.8
$ws.Cells.Item(87, 10).Value = 79750
$ws.Cells.Item(87, 11).Value = 86199.8
$ws.Cells.Item(87, 12).Value = 79750
$ws.Cells.Item(87, 13).Value = -84951.8
$ws.Cells.Item(87, 14).Value = -82246
# Row 90
$ws.Cells.Item(90, 8).Value = 84357
$ws.Cells.Item(90, 9).Value = 86199.8
$ws.Cells.Item(90, 10).Value = 79750
$ws.Cells.Item(90, 11).Value = 258599.4
$ws.Cells.Item(90, 12).Value = 239250
$ws.Cells.Item(90, 13).Value = -252359.4
$ws.Cells.Item(90, 14).Value = -251730

# ---------------------------------------------------------------------------
# CRP sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Cells.Item(58, 8).Value = 3256
$ws.Cells.Item(58, 9).Value = 2363.5715
$ws.Cells.Item(58, 11).Value = 2363.5715
$ws.Cells.Item(58, 13).Value = -2160.5715
# Row 86
$ws.Cells.Item(86, 8).Value = 4729.5
$ws.Cells.Item(86, 9).Value = 4259.8
$ws.Cells.Item(86, 11).Value = 4259.8
$ws.Cells.Item(86, 13).Value = -3136.8
# Row 89
$ws.Cells.Item(89, 8).Value = 4729.5
$ws.Cells.Item(89, 9).Value = 4259.8
$ws.Cells.Item(89, 11).Value = 21299
$ws.Cells.Item(89, 13).Value = -15683
# Row 104
$ws.Cells.Item(104, 8).Value = 54999.5
$ws.Cells.Item(104, 10).Value = 54999.5
$ws.Cells.Item(104, 12).Value = 54999.5
$ws.Cells.Item(104, 14).Value = -60241.5
# Row 105
$ws.Cells.Item(105, 8).Value = 10698.615
$ws.Cells.Item(105, 9).Value = 2134.3333
$ws.Cells.Item(105, 11).Value = 2134.3333
$ws.Cells.Item(105, 13).Value = -387.3332999999998
# Row 135
$ws.Cells.Item(135, 8).Value = 62467.8
$ws.Cells.Item(135, 10).Value = 62467.8
$ws.Cells.Item(135, 12).Value = 62467.8
$ws.Cells.Item(135, 14).Value = -72607.8
# Row 136
$ws.Cells.Item(136, 8).Value = 3256
$ws.Cells.Item(136, 9).Value = 2363.5715
$ws.Cells.Item(136, 11).Value = 7090.7145
$ws.Cells.Item(136, 13).Value = -4540.7145

# ---------------------------------------------------------------------------
# CUL sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 97
$ws.Cells.Item(97, 8).Value = 373.41666
$ws.Cells.Item(97, 10).Value = 497.66666
$ws.Cells.Item(97, 12).Value = 1492.99998
$ws.Cells.Item(97, 14).Value = -2484.99998
# Row 98
$ws.Cells.Item(98, 8).Value = 3398.75
$ws.Cells.Item(98, 9).Value = 1200
$ws.Cells.Item(98, 10).Value = 9995
$ws.Cells.Item(98, 11).Value = 3600
$ws.Cells.Item(98, 12).Value = 29985
$ws.Cells.Item(98, 13).Value = -2102
$ws.Cells.Item(98, 14).Value = -32981
# Row 104
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(104, 14).ClearContents()

# ---------------------------------------------------------------------------
# GSM sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 23
$ws.Cells.Item(23, 8).Value = 260.16666
$ws.Cells.Item(23, 9).Value = 171.5
$ws.Cells.Item(23, 10).Value = 437.5
$ws.Cells.Item(23, 11).Value = 171.5
$ws.Cells.Item(23, 12).Value = 437.5
$ws.Cells.Item(23, 13).Value = 51.5
$ws.Cells.Item(23, 14).Value = -883.5
# Row 70
$ws.Cells.Item(70, 8).Value = 6723.222
$ws.Cells.Item(70, 9).Value = 5500
$ws.Cells.Item(70, 11).Value = 5500
$ws.Cells.Item(70, 13).Value = -5230
# Row 73
$ws.Cells.Item(73, 8).Value = 6723.222
$ws.Cells.Item(73, 9).Value = 5500
$ws.Cells.Item(73, 11).Value = 5500
$ws.Cells.Item(73, 13).Value = -4564
# Row 80
$ws.Cells.Item(80, 8).Value = 3529.7
$ws.Cells.Item(80, 9).Value = 2649.5
$ws.Cells.Item(80, 10).Value = 4850
$ws.Cells.Item(80, 11).Value = 2649.5
$ws.Cells.Item(80, 12).Value = 4850
$ws.Cells.Item(80, 13).Value = -1651.5
$ws.Cells.Item(80, 14).Value = -6846
# Row 83
$ws.Cells.Item(83, 8).Value = 3529.7
$ws.Cells.Item(83, 9).Value = 2649.5
$ws.Cells.Item(83, 10).Value = 4850
$ws.Cells.Item(83, 11).Value = 13247.5
$ws.Cells.Item(83, 12).Value = 24250
$ws.Cells.Item(83, 13).Value = -8255.5
$ws.Cells.Item(83, 14).Value = -34234
# Row 95
$ws.Cells.Item(95, 8).Value = 413413
$ws.Cells.Item(95, 10).Value = 413413
$ws.Cells.Item(95, 12).Value = 413413
$ws.Cells.Item(95, 14).Value = -418905
# Row 122
$ws.Cells.Item(122, 8).Value = 6417.5386
$ws.Cells.Item(122, 9).Value = 3403.5
$ws.Cells.Item(122, 10).Value = 9001
$ws.Cells.Item(122, 11).Value = 10210.5
$ws.Cells.Item(122, 12).Value = 27003
$ws.Cells.Item(122, 13).Value = -7760.5
$ws.Cells.Item(122, 14).Value = -31903

# ---------------------------------------------------------------------------
# LTW sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Cells.Item(40, 8).Value = 3916.7144
$ws.Cells.Item(40, 9).Value = 3792.9375
$ws.Cells.Item(40, 11).Value = 3792.9375
$ws.Cells.Item(40, 13).Value = -3656.9375
# Row 46
$ws.Cells.Item(46, 8).Value = 1109.762
$ws.Cells.Item(46, 9).Value = 658.2
$ws.Cells.Item(46, 10).Value = 2238.6667
$ws.Cells.Item(46, 11).Value = 658.2
$ws.Cells.Item(46, 12).Value = 2238.6667
$ws.Cells.Item(46, 13).Value = -470.2
$ws.Cells.Item(46, 14).Value = -2614.6667
# Row 100
$ws.Cells.Item(100, 8).Value = 3194.9524
$ws.Cells.Item(100, 9).Value = 2792.8572
$ws.Cells.Item(100, 11).Value = 2792.8572
$ws.Cells.Item(100, 13).Value = -2251.8572
# Row 132
$ws.Cells.Item(132, 8).Value = 2857.24
$ws.Cells.Item(132, 9).Value = 2820.0454
$ws.Cells.Item(132, 10).Value = 3130
$ws.Cells.Item(132, 11).Value = 8460.1362
$ws.Cells.Item(132, 12).Value = 9390
$ws.Cells.Item(132, 13).Value = -5930.136200000001
$ws.Cells.Item(132, 14).Value = -14450
# Row 136
$ws.Cells.Item(136, 8).Value = 2325.5789
$ws.Cells.Item(136, 9).Value = 1762.3125
$ws.Cells.Item(136, 10).Value = 5329.6665
$ws.Cells.Item(136, 11).Value = 5286.9375
$ws.Cells.Item(136, 12).Value = 15988.9995
$ws.Cells.Item(136, 13).Value = -2736.9375
$ws.Cells.Item(136, 14).Value = -21088.9995

# ---------------------------------------------------------------------------
# WVR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 69
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 14).ClearContents()
# Row 72
$ws.Cells.Item(72, 8).Value = 0
$ws.Cells.Item(72, 10).Value = 0
$ws.Cells.Item(72, 12).Value = 0
$ws.Cells.Item(72, 14).ClearContents()
# Row 96
$ws.Cells.Item(96, 8).Value = 7499.5
$ws.Cells.Item(96, 10).Value = 10000
$ws.Cells.Item(96, 12).Value = 10000
$ws.Cells.Item(96, 14).Value = -12746
# Row 106
$ws.Cells.Item(106, 8).Value = 26700
$ws.Cells.Item(106, 9).Value = 23000
$ws.Cells.Item(106, 10).Value = 31633.334
$ws.Cells.Item(106, 11).Value = 23000
$ws.Cells.Item(106, 12).Value = 31633.334
$ws.Cells.Item(106, 13).Value = -21738
$ws.Cells.Item(106, 14).Value = -34157.334
# Row 136
$ws.Cells.Item(136, 8).Value = 1582.4062
$ws.Cells.Item(136, 9).Value = 1401.1852
$ws.Cells.Item(136, 11).Value = 4203.5556
$ws.Cells.Item(136, 13).Value = -1653.5556
